$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(19,19,16,17,18,20,18,20,16,22,26,21,19,18,16,21,22,13,22,23,20,17,15,18,16,21,17,17,20,14,21,20,24,21,26,22,14,20,19,20,18,9,15,17,13,24,16,21,19,13,21,15)
$colB = @(20,17,13,21,12,19,17,15,13,18,34,20,20,14,12,11,20,9,18,13,21,12,9,20,11,16,17,20,15,6,12,19,36,18,34,12,3,14,13,19,11,3,1,13,3,37,14,11,11,5,20,8)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

